$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selection on the sheet (cosmetic, matches diff's <selection activeCell="G11" sqref="G11"/>)
$ws.Range("G11").Select()

# Copy formatting from the D column cell in the description style used by D4
# (style index 21) onto D9, which previously used the empty-row style (15).
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in row 9 with new log entry data, matching the style/pattern of rows 7 and 8
$ws.Range("A9").Value = "Bugfix for cell array"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 45429
$ws.Range("D9").Value = "The cells were counted but not properly implemented with their actual`nwalls. This is now fixed."

# Match the row height used for the new content (diff shows ht="27")
$ws.Rows.Item(9).RowHeight = 27

$wb.Save()
